## Generate Report for Archive
## - Update localization status text "Ready for handoff" -> "In Translation"
##   (shared string is reused by Overview!E2:F3 and by the zh-cn/de-de
##   per-language "Status" column, cells C2:C3 on each language sheet).
## - The status column widths shrink to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: Status column (C) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("C3").Value = "In Translation"
$wsZh.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("C3").Value = "In Translation"
$wsDe.Columns.Item(3).ColumnWidth = 12.5

# --- Overview sheet: per-language status columns (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
